$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Bolinus brandaris" entry for the 2-RAP gear (old row 12) was removed;
# every subsequent data row in the 2-RAP block shifts up by one.
$ws.Rows("12:12").Delete()

# A handful of "Numb" (H) values that used to read 0 are now recorded as -1
# (Biological discard, Shells NA, Wood NA - now at rows 27, 42 and 45).
$ws.Range("H27").Value = -1
$ws.Range("H42").Value = -1
$ws.Range("H45").Value = -1

# The RF (I) column for every row in the 2-RAP block (rows 21:45) was
# recalculated against the new row count.
$ws.Range("I21:I45").Value = 39.12666666666667
